$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between the paired rows.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

function Swap-RowValues($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

# Swap data between row 3 and row 4
Swap-RowValues $ws 3 4 $cols

# Swap data between row 18 and row 19
Swap-RowValues $ws 18 19 $cols
